$d = $word.ActiveDocument

# 1) Remove the existing hidden "_GoBack" bookmark (it currently sits
#    further down in the document, inside the "modify ... and republish"
#    sentence near the end).
try {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
} catch {
    # no-op if it doesn't exist
}

# 2) Merge the paragraph that ends with " Version 1.1" with the following
#    (empty) paragraph, by deleting the paragraph mark between them.
$verPara = $d.Paragraphs.Item(2)
$verRange = $verPara.Range
$paraMark = $d.Range($verRange.End - 1, $verRange.End)
$paraMark.Delete()

# 3) Re-add the "_GoBack" bookmark right after the "Version 1.1" text,
#    i.e. at the end of the (now merged) paragraph's text.
$verPara2 = $d.Paragraphs.Item(2)
$verRange2 = $verPara2.Range
$insertPoint = $d.Range($verRange2.End - 1, $verRange2.End - 1)
$d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null
